# Update the "Integer min" value for rule R30 (row 10) in the Rules sheet
# from 18 to 1, per the authored revision restore.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
